# Incomplete 019 - Redo
# Add a new (incomplete) entry for the "lengthOfLongestSubstring" problem
# (day 13, revisited/redone) to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the Description first, then day / RoadMap? / Type / Name - this is
# the order the cells were actually typed in, and it controls the order new
# entries land in the shared-strings table.
$ws.Range("F23").Value = "Typical exercise for sliding window algorithm, increases the window while letter is different, reset once it is the same "
$ws.Range("B23").Value = "13 missed 1"
$ws.Range("C23").Value = "REDO"
$ws.Range("D23").Value = "Sliding window"
$ws.Range("E23").Value = "lengthOfLongestSubstring"

# Highlight the "Name" cell the same way other entries in this column are highlighted
$ws.Range("E23").Interior.Color = $ws.Range("E22").Interior.Color

# Reflect the wider "day" column needed to fit the new text, and move the
# active selection to where the user left off editing
$ws.Columns("B").ColumnWidth = 11.1667
$ws.Range("F28").Select()
